$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2  = @(0.230598549219447, 0.7445381377257259, 1.018876335886205, 1.009394043912587, 1.011190667123304, 18)
  3  = @(0.3045708550604296, 0.3907590825865097, 0.2095424611763801, 0.4577580814976182, 0.3522467123879903, 17)
  4  = @(0.2502419823634077, 0.366077697295096, 0.1747842588224781, 0.4180720737175327, 0.345891249322188, 16)
  5  = @(0.3395964577034516, 0.406420661795407, 0.2182909649593746, 0.4672161865340012, 0.3321443493336289, 15)
  6  = @(0.3836894047943985, 0.3861383846528136, 0.2113509978471021, 0.4597292658153297, 0.2628055933924591, 14)
  7  = @(0.3925847413104661, 0.4061565136574045, 0.2245340717459249, 0.4738502630007974, 0.2761863507775066, 13)
  8  = @(0.4359548671747835, 0.4359548671747835, 0.2345424034301451, 0.4842957809336615, 0.2202950679514644, 12)
  9  = @(0.4182767482978756, 0.4182767482978756, 0.2153166270150599, 0.4640222268545548, 0.2107066864939436, 11)
  10 = @(0.3824308908543818, 0.3824308908543818, 0.1807831973670264, 0.4251860738159546, 0.1958735734184745, 10)
  11 = @(0.3843625910566515, 0.3843625910566515, 0.1817866140615611, 0.4263644146285676, 0.195725609566047, 9)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
    $ws.Cells.Item($r, 7).Value = $vals[5]
}
